$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- Edit sheet: 展览 (sheet1) ----
$ws1.Rows.Item(27).Insert()

# New row 27: 九江·SXD动漫嘉年华
$ws1.Range("A27").Value = 26
$ws1.Range("B27").Value = "'2024-07-21"
$ws1.Range("C27").Value = "九江·SXD动漫嘉年华"
$ws1.Range("D27").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws1.Range("E27").Value = "2024.07.21 10:00-07.21 17:30"
$ws1.Range("F27").Value = 4
$ws1.Range("G27").Value = 45
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws1.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

# Simple F-column increments (rows 2-25, unaffected by the row insert)
$ws1.Range("F2").Value = 196
$ws1.Range("F3").Value = 3149
$ws1.Range("F4").Value = 228
$ws1.Range("F5").Value = 121
$ws1.Range("F7").Value = 1668
$ws1.Range("F8").Value = 1625
$ws1.Range("F13").Value = 190
$ws1.Range("F18").Value = 5
$ws1.Range("F24").Value = 188
$ws1.Range("F25").Value = 99
$ws1.Range("F9").Value = 344

# Text change carried with the shifted row (old row38 -> new row39)
$ws1.Range("C39").Value = "宜春·第三十五届静卿国风动漫文化展览会"

# F-column tweaks in the shifted tail (rows 34,36,39,40,41 post-shift)
$ws1.Range("F34").Value = 208
$ws1.Range("F36").Value = 303
$ws1.Range("F39").Value = 296
$ws1.Range("F40").Value = 507
$ws1.Range("F41").Value = 288

# ---- Edit sheet: 全部类型 (sheet4) ----
$ws4.Rows.Item(27).Insert()

# New row 27: 九江·SXD动漫嘉年华
$ws4.Range("A27").Value = 26
$ws4.Range("B27").Value = "'2024-07-21"
$ws4.Range("C27").Value = "九江·SXD动漫嘉年华"
$ws4.Range("D27").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws4.Range("E27").Value = "2024.07.21 10:00-07.21 17:30"
$ws4.Range("F27").Value = 4
$ws4.Range("G27").Value = 45
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws4.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

# Simple F-column increments (rows 2-25, unaffected by the row insert)
$ws4.Range("F2").Value = 196
$ws4.Range("F3").Value = 3149
$ws4.Range("F4").Value = 228
$ws4.Range("F5").Value = 121
$ws4.Range("F7").Value = 1668
$ws4.Range("F8").Value = 1625
$ws4.Range("F13").Value = 190
$ws4.Range("F18").Value = 5
$ws4.Range("F24").Value = 188
$ws4.Range("F25").Value = 99
$ws4.Range("F9").Value = 346

# Text change carried with the shifted row (old row38 -> new row39)
$ws4.Range("C39").Value = "宜春·第三十五届静卿国风动漫文化展览会"

# F-column tweaks in the shifted tail (rows 34,36,39,40,41 post-shift)
$ws4.Range("F34").Value = 208
$ws4.Range("F36").Value = 305
$ws4.Range("F39").Value = 299
$ws4.Range("F40").Value = 507
$ws4.Range("F41").Value = 290

Write-Output "done"
